$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix typo in e011a cell (B21): add trailing space after "Die Roll ="
$b21Text = @'
<Bold>e011a Deployment - Counterattack Scenario</Bold> 
<InlineUIContainer><Button Content='r20.41' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Since this is a Counterattack scenario, your tank&apos;s deployment is automatically stopped and hull down.  You must still roll for if your tank is the lead tank per the 
<InlineUIContainer><Button Content='Deployment' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  Table. 
<LineBreak/><LineBreak/>
Die Roll = 
<InlineUIContainer><Image Name='DieRollBlue' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$ws.Range("B21").Value = $b21Text

# 2) Add Resistance/Area Type lines to e032a cell (B45)
$b45Text = @'
<Bold>e032a Battle Check - Counterattack</Bold> 
<InlineUIContainer><Button Content='r20.42' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
Time: TIME_OF_DAY   
Resistance: RESISTANCE_OF_DAY<LineBreak/>
Area Type: AREA_TYPE<LineBreak/>
 <LineBreak/><LineBreak/>
Choose 
<InlineUIContainer><Button Content='Resupply' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> or 
roll 1D and consult the <InlineUIContainer><Button Content='Resistance' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table to determine if combat occurs by counterattacking German forces: 
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
'@
$ws.Range("B45").Value = $b45Text

# Row 45 auto-grows from the extra lines of text (150 -> 180)
$ws.Rows.Item(45).RowHeight = 180

# 3) Remove obsolete e033a row (old row 47); this shifts subsequent rows up
$ws.Rows.Item(47).Delete()
